# Applies the "AutreActiviteSoumiseAutorisation.EntiteGeographique" row addition
# plus the metadata Date bump described by the commit diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Metadata sheet: bump the "Date" value (row 8, column B)
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# ---------------------------------------------------------------------------
# 2) Elements sheet: append a new row (row 10) describing the
#    "EntiteGeographique" reference element, cloned from the formatting of
#    the last existing data row (row 9) so styles/column types line up.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Clone the style of row 9 into row 10 (keeps cell style s="2" without
# minting any new style records).
$ws.Range("A9:AJ9").Copy()
$ws.Range("A10:AJ10").PasteSpecial(-4122)

# ID / Path / Base Path -> new element path text
$id = "AutreActiviteSoumiseAutorisation.EntiteGeographique"
$ws.Cells.Item(10, 1).Value = $id
$ws.Cells.Item(10, 2).Value = $id
$ws.Cells.Item(10, 32).Value = $id

# Min / Max / Base Min / Base Max -> text "1" (copied from an existing cell
# that already stores "1" as text so we reuse the shared string & avoid
# Excel re-typing it as a number).
$ws.Range("G3").Copy()
$ws.Cells.Item(10, 6).PasteSpecial(-4163)
$ws.Cells.Item(10, 7).PasteSpecial(-4163)
$ws.Cells.Item(10, 33).PasteSpecial(-4163)
$ws.Cells.Item(10, 34).PasteSpecial(-4163)

# Type(s) -> link to the EntiteGeographique StructureDefinition
$ws.Cells.Item(10, 11).Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/EntiteGeographique`n"

# Short / Definition -> description text
$desc = "Lien vers la classe EntiteGeographique"
$ws.Cells.Item(10, 12).Value = $desc
$ws.Cells.Item(10, 13).Value = $desc

# Widen the "Type(s)" column (K) so the new long URL value fits, matching
# the authored column width.
$ws.Columns.Item(11).ColumnWidth = 60.83

Write-Output "AutreActiviteSoumiseAutorisation.EntiteGeographique row added"
